$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F551").Value = 18246
$ws.Range("F558").Value = 24841
$ws.Range("F575").Value = 26443
$ws.Range("F614").Value = 47744
$ws.Range("F616").Value = 38194
$ws.Range("G616").Value = 2570
$ws.Range("F618").Value = 37792
$ws.Range("G618").Value = 2658
$ws.Range("F621").Value = 55980
$ws.Range("G621").Value = 4120
$ws.Range("F622").Value = 41315
$ws.Range("G622").Value = 3010
$ws.Range("F623").Value = 15009
$ws.Range("G623").Value = 1573
$ws.Range("F624").Value = 51061
$ws.Range("G624").Value = 3938
$ws.Range("F625").Value = 43674
$ws.Range("G625").Value = 3558
$ws.Range("F626").Value = 20016
$ws.Range("F627").Value = 33588
$ws.Range("G627").Value = 2719
$ws.Range("F628").Value = 64210
$ws.Range("G628").Value = 4180
$ws.Range("F629").Value = 46045
$ws.Range("G629").Value = 2916
$ws.Range("F630").Value = 46436
$ws.Range("G630").Value = 2944
$ws.Range("F631").Value = 41694
$ws.Range("G631").Value = 2776
$ws.Range("F632").Value = 43937
$ws.Range("G632").Value = 2628
$ws.Range("F633").Value = 23868
$ws.Range("G633").Value = 1920
$ws.Range("F634").Value = 45981
$ws.Range("G634").Value = 2148
$ws.Range("F635").Value = 82014
$ws.Range("G635").Value = 3646
$ws.Range("F636").Value = 48912
$ws.Range("G636").Value = 2293
$ws.Range("F637").Value = 42794
$ws.Range("G637").Value = 2047
$ws.Range("F638").Value = 36486
$ws.Range("G638").Value = 2283
$ws.Range("F639").Value = 34882
$ws.Range("G639").Value = 1752
$ws.Range("F640").Value = 18036
$ws.Range("G640").Value = 1168
$ws.Range("F641").Value = 25716
$ws.Range("G641").Value = 1077
